$wb = $excel.ActiveWorkbook

# The change needs to be applied to both the "展览" and "全部类型" sheets,
# which contain identical data in this workbook.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1666
    $ws.Range("F6").Value = 441
    $ws.Range("F9").Value = 564
}
